# Normalize phone number formatting in the "Phone" column (D) of the
# Legal Advocate services sheet. The raw/inconsistent phone strings are
# replaced with a consistent "(NNN) NNN-NNNN"-style format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "(416) 268-2318"

$ws.Range("D14").Value = "(888) 232-0232"
$ws.Range("D15").Value = "(888) 232-0232"
$ws.Range("D16").Value = "(888) 232-0232"
$ws.Range("D18").Value = "(888) 232-0232"
$ws.Range("D19").Value = "(888) 232-0232"
$ws.Range("D20").Value = "(888) 232-0232"
$ws.Range("D21").Value = "(888) 232-0232"
$ws.Range("D23").Value = "(888) 232-0232"
$ws.Range("D24").Value = "(888) 232-0232"

$ws.Range("D30").Value = "(709) 722-4031"

$ws.Range("D32").Value = "(888) 508-3028-101"

$ws.Range("D43").Value = "(866) 845-3425"
$ws.Range("D44").Value = "(866) 845-3425"

$ws.Range("D47").Value = "(867) 979-2228"
$ws.Range("D48").Value = "(250) 635-3178"
$ws.Range("D49").Value = "(905) 689-4727"
$ws.Range("D50").Value = "(519) 752-5308"
